# Add three new "KPI" textboxes to slide 1, right after the existing
# trailing "#" placeholder textboxes, matching the authored template.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- KPI1 -----------------------------------------------------------
$kpi1 = $s.Shapes.AddTextbox(1, 2.274015748031496, 47.164724409448816, 181.7003937007874, 29.081259842519685)
$kpi1.TextFrame.WordWrap = -1
$kpi1.TextFrame.AutoSize = 1
$kpi1.Fill.Visible = 0

$kpi1.TextFrame.TextRange.Text = "KPI1"
$kpi1.TextFrame.TextRange.ParagraphFormat.Alignment = 1

$kpi1.TextFrame.TextRange.Font.Size = 18
$kpi1.TextFrame.TextRange.Font.Bold = 0
$kpi1.TextFrame.TextRange.Font.Italic = 0
$kpi1.TextFrame.TextRange.Font.Underline = 0
$kpi1.TextFrame.TextRange.Font.Shadow = 0
$kpi1.TextFrame.TextRange.Font.Color.RGB = 0xE6E6E7

$kpi1f2 = $kpi1.TextFrame2.TextRange.Font
$kpi1f2.Name = "Work Sans ExtraLight"
$kpi1f2.NameFarEast = "+mn-ea"
$kpi1f2.NameComplexScript = "+mn-cs"

# --- KPI2 -----------------------------------------------------------
$kpi2 = $s.Shapes.AddTextbox(1, 2.990472440944882, 189.96220472440945, 180.26740157480316, 29.081259842519685)
$kpi2.TextFrame.WordWrap = -1
$kpi2.TextFrame.AutoSize = 1
$kpi2.Fill.Visible = 0

$kpi2.TextFrame.TextRange.Text = "KPI2"
$kpi2.TextFrame.TextRange.IndentLevel = 1

$kpi2.TextFrame.TextRange.Font.Name = "Work Sans ExtraLight"
$kpi2.TextFrame.TextRange.Font.Color.RGB = 0xE6E6E7

# --- KPI3 -----------------------------------------------------------
$kpi3 = $s.Shapes.AddTextbox(1, 5.940629921259842, 342.49952755905514, 177.3172440944882, 29.081259842519685)
$kpi3.TextFrame.WordWrap = -1
$kpi3.TextFrame.AutoSize = 1
$kpi3.Fill.Visible = 0

$kpi3.TextFrame.TextRange.Text = "KPI3"

$kpi3.TextFrame.TextRange.Font.Name = "Work Sans ExtraLight"
$kpi3.TextFrame.TextRange.Font.Color.RGB = 0xE6E6E7
